$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the H1
#    title paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$metaEnd   = $metaPara.Range.End
$metaFull  = $d.Range($metaStart, $metaEnd)
$metaFull.InsertBefore("Meta description: Try Candy Witch for free and discover two exciting bonus games with Sticky Wilds and progressive multipliers, plus great payouts.")

# Bold just the "Meta description" label (16 characters) -- leave the
# rest (starting at the colon) unformatted.
$labelStart = $metaPara.Range.Start
$labelRange = $d.Range($labelStart, $labelStart + 16)
$labelRange.Bold = 1

# ------------------------------------------------------------------
# 2. The document used to end with two paragraphs duplicating the
#    title (bold) and the meta blurb (italic). Drop the stray bold
#    title paragraph entirely, and turn the italic blurb into the new
#    image-generation prompt.
# ------------------------------------------------------------------
$countBeforeTrim = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($countBeforeTrim - 1)
$dupTitlePara.Range.Delete()

$blurbPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$blurbFind = $blurbPara.Range.Find
$blurbFind.Execute(
    "Try Candy Witch for free and discover two exciting bonus games with Sticky Wilds and progressive multipliers, plus great payouts.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create a cartoon-style feature image for Candy Witch that features a happy Maya warrior with glasses. The image should have a bright and colorful background, with the Maya warrior positioned in the center of the frame. He should be smiling and holding a handful of candy in one hand, with the other hand raised up in a magic spell-casting pose. His glasses should be oversized and cartoonish, with a reflection of the glow from the enchanted forest in the lenses. In the background, there should be hints of the forest and the moonlight. The overall image should convey a sense of fun and whimsy while also capturing the magic and excitement of the Candy Witch slot game.",
    2) | Out-Null
